# Weekly Fruta/Hortaliza update: two new price records were added for
# "Comercializadora del Agro de Limarí - Pepino ensalada", inserted right
# after the existing row for the Coquimbo region (row 66), which pushes
# all subsequent rows down by two positions (old row 67 -> new row 69,
# ..., old row 92 -> new row 94). The sheet's used range grows from
# A1:R92 to A1:R94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the old row 67, shifting everything
# below it (rows 67-92) down to rows 69-94.
$ws.Rows("67:68").Insert()

# New row 67 - "Primera" quality record dated 2021-11-24 (serial 44524)
$ws.Cells.Item(67,1).Value = 2
$ws.Cells.Item(67,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(67,3).Value = "Coquimbo"
$ws.Cells.Item(67,4).Value = 44524
$ws.Cells.Item(67,5).Value = 4
$ws.Cells.Item(67,6).Value = 100112043
$ws.Cells.Item(67,7).Value = "Pepino ensalada"
$ws.Cells.Item(67,8).Value = "Sin especificar"
$ws.Cells.Item(67,9).Value = "Primera"
$ws.Cells.Item(67,10).Value = 700
$ws.Cells.Item(67,11).Value = 6500
$ws.Cells.Item(67,12).Value = 7000
$ws.Cells.Item(67,13).Value = 6750
$ws.Cells.Item(67,14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(67,15).Value = "Provincia de Limarí"
$ws.Cells.Item(67,16).Value = 96
$ws.Cells.Item(67,17).Value = 70
$ws.Cells.Item(67,18).Value = "Hortaliza"

# New row 68 - "Segunda" quality record dated 2021-11-24 (serial 44524)
$ws.Cells.Item(68,1).Value = 2
$ws.Cells.Item(68,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(68,3).Value = "Coquimbo"
$ws.Cells.Item(68,4).Value = 44524
$ws.Cells.Item(68,5).Value = 4
$ws.Cells.Item(68,6).Value = 100112043
$ws.Cells.Item(68,7).Value = "Pepino ensalada"
$ws.Cells.Item(68,8).Value = "Sin especificar"
$ws.Cells.Item(68,9).Value = "Segunda"
$ws.Cells.Item(68,10).Value = 400
$ws.Cells.Item(68,11).Value = 4500
$ws.Cells.Item(68,12).Value = 5000
$ws.Cells.Item(68,13).Value = 4750
$ws.Cells.Item(68,14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(68,15).Value = "Provincia de Limarí"
$ws.Cells.Item(68,16).Value = 48
$ws.Cells.Item(68,17).Value = 100
$ws.Cells.Item(68,18).Value = "Hortaliza"
